$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: add label for the integration-by-parts function in C12 (red font, like G15/H15 below)
$ws.Range("C12").Value = "4*e^(-2*t)"
$ws.Range("C12").Font.Color = 255
$ws.Range("C12").Font.Name = "Aptos"

# Row 15: fix the formula that used to error (#NUM!) and add the evaluated-bounds labels
$ws.Range("C15").Formula = "=SQRT(10+6)"

$ws.Range("G15").Value = "e^2*(1)"
$ws.Range("G15").Font.Color = 255
$ws.Range("G15").Font.Name = "Aptos"

$ws.Range("H15").Value = "e^0*(-1)"
$ws.Range("H15").Font.Color = 255
$ws.Range("H15").Font.Name = "Aptos"

# Row 16: evaluate the two bound terms of the integral
$ws.Range("G16").Formula = "=EXP(2)"
$ws.Range("H16").Formula = "=EXP(0)*(-1)"

# Row 17: subtract to get the final definite integral result
$ws.Range("G17").Formula = "=G16-H16"

# Update the active selection to reflect where editing ended
$ws.Range("G17").Select() | Out-Null
